# ---------------------------------------------------------------------------
# AAPLBag.xlsx update
#
# The trading-data export routine changed: the sentiment score that used to
# be saved as a single "totalScore" column is now split into a final blended
# "ScoreFinal" (technical + fundamental + sentiment) and the raw
# "totalSentiment" value, and two new fundamental indicators ("RSI", "PEG")
# are appended after "Method". The exported data itself is refreshed with a
# new run that wrote two rows (a half trade and a normal trade) instead of
# the old single row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----------------------------------------------------
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "ScoreFinal"
$ws.Range("C1").Value = "totalSentiment"
$ws.Range("D1").Value = "posWordPercentage"
$ws.Range("E1").Value = "negWordPercentage"
$ws.Range("F1").Value = "posPhrasePercentage"
$ws.Range("G1").Value = "negPhrasePercentage"
$ws.Range("H1").Value = "ElapsedMs"
$ws.Range("I1").Value = "wordCount"
$ws.Range("J1").Value = "sentenceCount"
$ws.Range("K1").Value = "posWordCount"
$ws.Range("L1").Value = "negWordCount"
$ws.Range("M1").Value = "positivePhraseCount"
$ws.Range("N1").Value = "negativePhraseCount"
$ws.Range("O1").Value = "Method"
$ws.Range("P1").Value = "RSI"
$ws.Range("Q1").Value = "PEG"

# ---- Data row 2 -------------------------------------------------------------
$ws.Range("A2").Value = 42627.87427083333
$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 36
$ws.Range("D2").Value = 64
$ws.Range("E2").Value = 33
$ws.Range("F2").Value = 99
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 5231
$ws.Range("I2").Value = 7031
$ws.Range("J2").Value = 743
$ws.Range("K2").Value = 122
$ws.Range("L2").Value = 64
$ws.Range("M2").Value = 12
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = "Bag"
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 1

# ---- Data row 3 (new row, second saved trade of this run) -------------------
$ws.Range("A3").Value = 42627.877280092594
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 63
$ws.Range("E3").Value = 35
$ws.Range("F3").Value = 100
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 3067
$ws.Range("I3").Value = 5351
$ws.Range("J3").Value = 545
$ws.Range("K3").Value = 81
$ws.Range("L3").Value = 46
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = "Bag"
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 1

# ---- Column widths (best-fit, follows the new header text) ------------------
# The runtime quantizes ColumnWidth assignments onto a 1/6-character grid
# before it is written back out as the OOXML `width` attribute, so the
# values below are chosen as the closest achievable pre-images of the real
# (sub-pixel) best-fit widths produced by a real Excel session.
$ws.Columns("A:A").ColumnWidth = 14
$ws.Columns("B:B").ColumnWidth = 9.333333333333332
$ws.Columns("C:C").ColumnWidth = 13.666666666666666
$ws.Columns("D:D").ColumnWidth = 18.5
$ws.Columns("E:E").ColumnWidth = 18.666666666666668
$ws.Columns("F:F").ColumnWidth = 19.666666666666668
$ws.Columns("G:G").ColumnWidth = 19.666666666666668
$ws.Columns("H:H").ColumnWidth = 9.666666666666666
$ws.Columns("I:I").ColumnWidth = 10
$ws.Columns("J:J").ColumnWidth = 13.666666666666666
$ws.Columns("K:K").ColumnWidth = 13.666666666666666
$ws.Columns("L:L").ColumnWidth = 13.666666666666666
$ws.Columns("M:M").ColumnWidth = 18.833333333333336

Write-Output "AAPLBag.xlsx updated"
